# "Fruta / hortaliza, semanal" update
#
# A new daily price record (date serial 45106 = 2023-06-29) is added for the
# "Plátano" product block on this sheet. Because the underlying data feed
# keeps each "Pintón" / "Primera Pintón" quality pair together, the new
# observation is inserted at the TOP of that block (rows 999:1000), pushing
# every subsequent row down by two. The tail of the block (previously rows
# 1084:1085) ends up duplicated at the new end of the table (rows 1086:1087).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new observations - shifts old rows 999-1085 down to
# 1001-1087 (and grows the sheet's used range to A1:T1087 accordingly).
$ws.Rows("999:1000").Insert()

# New row 999 - "Pintón" quality, same market/product metadata as its
# neighbours, price data unchanged from what used to be here, just a newer
# date.
$ws.Range("A999").Value = 5
$ws.Range("B999").Value = "Macroferia Regional de Talca"
$ws.Range("C999").Value = "Maule"
$ws.Range("D999").Value = 45106
$ws.Range("E999").Value = 7
$ws.Range("F999").Value = "Fruta"
$ws.Range("G999").Value = 100108
$ws.Range("H999").Value = "Tropicales y subtropicales"
$ws.Range("I999").Value = 100108006
$ws.Range("J999").Value = "Plátano"
$ws.Range("K999").Value = "Sin especificar"
$ws.Range("L999").Value = "Pintón"
$ws.Range("M999").Value = 800
$ws.Range("N999").Value = 11000
$ws.Range("O999").Value = 11000
$ws.Range("P999").Value = 11000
$ws.Range("Q999").Value = "$/caja 20 kilos"
$ws.Range("R999").Value = "Ecuador"
$ws.Range("S999").Value = 550
$ws.Range("T999").Value = 20

# New row 1000 - "Primera Pintón" quality, paired with row 999, with an
# updated volume (M) for the new date.
$ws.Range("A1000").Value = 5
$ws.Range("B1000").Value = "Macroferia Regional de Talca"
$ws.Range("C1000").Value = "Maule"
$ws.Range("D1000").Value = 45106
$ws.Range("E1000").Value = 7
$ws.Range("F1000").Value = "Fruta"
$ws.Range("G1000").Value = 100108
$ws.Range("H1000").Value = "Tropicales y subtropicales"
$ws.Range("I1000").Value = 100108006
$ws.Range("J1000").Value = "Plátano"
$ws.Range("K1000").Value = "Sin especificar"
$ws.Range("L1000").Value = "Primera Pintón"
$ws.Range("M1000").Value = 600
$ws.Range("N1000").Value = 12000
$ws.Range("O1000").Value = 12000
$ws.Range("P1000").Value = 12000
$ws.Range("Q1000").Value = "$/caja 20 kilos"
$ws.Range("R1000").Value = "Ecuador"
$ws.Range("S1000").Value = 600
$ws.Range("T1000").Value = 20
